# Bill of Materials update:
#  - Row 9  (PCB spare part "26"): amount 1 -> 0 (no longer needed)
#  - Row 16 (item "43"): amount 0 -> 1 (now needed)
#  - Row 17 (item "45"): amount 0 -> 1 (now needed)
#  - Row 26 (item "24"/Mosfet-ish row referencing string 70): amount blank -> 1 (now needed)
# The F column formulas (=Dn*En) and the F30 total (=SUM(F3:F29)) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("D26").Value = 1

# Reflect the author's final selection on the totals cell.
$ws.Range("F30").Select()
